$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 carries the default (unstyled) cell format in this sheet. We use it
# as a template to strip the temporary "@" text format back off any cell
# whose new numeric-looking value would otherwise be auto-converted by
# Excel into a real number (losing formatting like trailing zeros or the
# multi-dot "thousands" separators used in this sheet).
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "66.713.87"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "3.256.40"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.96"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.21"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.255.88"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.548"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.93"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +7.41%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.32"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "3.790.93"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "66.754.11"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "3.259.06"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "507.72"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.42"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("E22").Value = "  +3.30%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.81"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.68"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.04"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +56.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.09"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.90"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -6.02%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.13"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  +21.78%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0796"
$ws.Range("E38").Value = "  +18.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.81"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "495.48"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").Value = "2.968.62"
$ws.Range("E46").Value = "  +4.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.76"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.48"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  +4.69%  "
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.60"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -0.37%  "
